$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.550712704658508
$ws.Range("B1").Value = 4.22331428527832
$ws.Range("C1").Value = 3.34261155128479
$ws.Range("D1").Value = 1.929315447807312
$ws.Range("E1").Value = 0.7013711929321289
